# Update odds values on Sheet1 per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 2.5

# Row 4
$ws.Range("T4").Value = 2

# Row 5
$ws.Range("G5").Value = 1.7
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 5.5
$ws.Range("J5").Value = 2.4
$ws.Range("U5").Value = 2.5
$ws.Range("V5").Value = 1.5
$ws.Range("Z5").Value = 12
$ws.Range("AA5").Value = 17
$ws.Range("AI5").Value = 26
$ws.Range("AK5").Value = 67
$ws.Range("AN5").Value = 3.4
$ws.Range("AW5").Value = 7
$ws.Range("AZ5").Value = 151
